# "Generate Report for handback"
#
# For each localized-language sheet (zh-cn, de-de) this:
#   - updates the status text from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two real file rows,
#   - fills in the "Latest Target File" (E) / "Latest Handback File" (F)
#     columns (mirroring the source file / handoff file respectively),
#     including hyperlinks that point at the same targets as columns A/C,
#   - stamps the "Latest Handback DateTime" (G) column with the actual
#     handback time (previously the 0001-01-01 00:00:00 placeholder).
# The Overview sheet pulls the same shared status text, so updating it
# there too keeps everything in sync.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# --- 1. Overview sheet: refresh the status text in both language columns ---
$overview = $wb.Worksheets.Item("Overview")
if ($overview.Range("B2").Text -eq $oldStatus) { $overview.Range("B2").Value = $newStatus }
if ($overview.Range("C2").Text -eq $oldStatus) { $overview.Range("C2").Value = $newStatus }
if ($overview.Range("B3").Text -eq $oldStatus) { $overview.Range("B3").Value = $newStatus }
if ($overview.Range("C3").Text -eq $oldStatus) { $overview.Range("C3").Value = $newStatus }

# Helper: look up the URL a hyperlink on a given cell points to.
function Get-HyperlinkAddress($sheet, $row, $col) {
    foreach ($hl in $sheet.Hyperlinks) {
        if (($hl.Range.Row -eq $row) -and ($hl.Range.Column -eq $col)) {
            return $hl.Address
        }
    }
    return $null
}

# Per-sheet handback timestamps recorded for this report.
$handbackTimes = @{
    "zh-cn" = "2016-01-22 03:07:52"
    "de-de" = "2016-01-22 03:08:16"
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $handbackTime = $handbackTimes[$sheetName]

    # Capture the existing (A/C/A4) hyperlink targets before we wipe the
    # collection, then rebuild the whole thing -- including the new E/F
    # links -- in row-major, left-to-right order, same as a freshly
    # generated report would lay them out.
    $addrA2 = Get-HyperlinkAddress $ws 2 1
    $addrC2 = Get-HyperlinkAddress $ws 2 3
    $addrA3 = Get-HyperlinkAddress $ws 3 1
    $addrC3 = Get-HyperlinkAddress $ws 3 3
    $addrA4 = Get-HyperlinkAddress $ws 4 1

    $dispA2 = $ws.Cells.Item(2, 1).Text
    $dispC2 = $ws.Cells.Item(2, 3).Text
    $dispA3 = $ws.Cells.Item(3, 1).Text
    $dispC3 = $ws.Cells.Item(3, 3).Text
    $dispA4 = $ws.Cells.Item(4, 1).Text

    foreach ($row in 2, 3) {
        # Status column
        if ($ws.Cells.Item($row, 2).Text -eq $oldStatus) {
            $ws.Cells.Item($row, 2).Value = $newStatus
        }

        $sourceFile = $ws.Cells.Item($row, 1).Text
        $handoffFile = $ws.Cells.Item($row, 3).Text

        # E: Latest Target File -- mirrors the source file (column A)
        $ws.Cells.Item($row, 5).Value = $sourceFile

        # F: Latest Handback File -- mirrors the handoff file (column C)
        $ws.Cells.Item($row, 6).Value = $handoffFile

        # G: Latest Handback DateTime -- now a real timestamp
        $ws.Cells.Item($row, 7).Value = $handbackTime
    }

    $ws.Hyperlinks.Delete()
    if ($addrA2 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(2, 1), $addrA2, "", "", $dispA2) }
    if ($addrC2 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(2, 3), $addrC2, "", "", $dispC2) }
    if ($addrA2 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(2, 5), $addrA2, "", "", $dispA2) }
    if ($addrC2 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(2, 6), $addrC2, "", "", $dispC2) }
    if ($addrA3 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(3, 1), $addrA3, "", "", $dispA3) }
    if ($addrC3 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(3, 3), $addrC3, "", "", $dispC3) }
    if ($addrA3 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(3, 5), $addrA3, "", "", $dispA3) }
    if ($addrC3 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(3, 6), $addrC3, "", "", $dispC3) }
    if ($addrA4 -ne $null) { $ws.Hyperlinks.Add($ws.Cells.Item(4, 1), $addrA4, "", "", $dispA4) }
}
